$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Jugadores")
$ws2 = $wb.Worksheets.Item("Entrenadores")

# --- Sheet "Entrenadores": insert a new column before I ---
[void]$ws2.Columns("I").Insert()

# New header cell (shares the header style used by other header cells, e.g. H1:
# centered horizontally/vertically)
$ws2.Range("I1").Value = "Nombre Foto Plantel Club"
$ws2.Range("I1").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("I1").VerticalAlignment = -4108    # xlCenter

# Column I width (match column H's width, like the rest of the header columns)
$ws2.Columns("I").ColumnWidth = $ws2.Columns("H").ColumnWidth

# --- View/selection bookkeeping ---
# Scroll sheet "Jugadores" so column AJ is the left-most visible column,
# while keeping its existing selection (K2) untouched.
[void]$ws1.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 36

# Make "Entrenadores" the active sheet/tab and select I1 on it.
[void]$ws2.Activate()
[void]$ws2.Range("I1").Select()

Write-Host "done"
